$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9999994693310038
$ws.Range("E2").Value = 0.9999994693310038

# Row 3
$ws.Range("D3").Value = 0.9999996710763952
$ws.Range("E3").Value = 0.9999996710763952

# Row 4
$ws.Range("D4").Value = 0.3994707376023923
$ws.Range("E4").Value = 0.3994707376023923

# Row 5
$ws.Range("D5").Value = 0.966933350705273
$ws.Range("E5").Value = 0.966933350705273

# Row 6
$ws.Range("D6").Value = 0.985380069886631
$ws.Range("E6").Value = 0.985380069886631

# Row 7 (C7 unchanged, true)
$ws.Range("D7").Value = 0.8699587349927055
$ws.Range("E7").Value = 0.1300412650072945

# Row 8 (C8 flips true -> false)
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = [double]"2.416037725771704E-07"
$ws.Range("E8").Value = 0.9999997583962275

# Row 9 (C9 flips true -> false)
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = [double]"6.096162340735897E-06"
$ws.Range("E9").Value = 0.9999939038376593

# Row 10 (C10 flips true -> false)
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = [double]"2.811749717891434E-07"
$ws.Range("E10").Value = 0.9999997188250282

# Row 11 (C11 flips true -> false)
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.3864789523720936
$ws.Range("E11").Value = 0.6135210476279064
$ws.Range("F11").Value = 8.09391975402832
$ws.Range("G11").Value = 0.2

# Row 12
$ws.Range("D12").Value = 0.9999999857401505
$ws.Range("E12").Value = 0.9999999857401505

# Row 13
$ws.Range("D13").Value = 0.9999999999954774
$ws.Range("E13").Value = 0.9999999999954774

# Row 14
$ws.Range("D14").Value = 0.710221910929034
$ws.Range("E14").Value = 0.710221910929034

# Row 15
$ws.Range("D15").Value = [double]"8.051458931515577E-05"
$ws.Range("E15").Value = [double]"8.051458931515577E-05"

# Row 16
$ws.Range("D16").Value = 0.9984148830855072
$ws.Range("E16").Value = 0.9984148830855072

# Row 17 (C17 unchanged, true)
$ws.Range("D17").Value = 0.9804361248343552
$ws.Range("E17").Value = 0.01956387516564484

# Row 18 (C18 flips true -> false)
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = [double]"2.501861133084964E-10"
$ws.Range("E18").Value = 0.9999999997498139

# Row 19 (C19 flips true -> false)
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = [double]"2.290886819987316E-06"
$ws.Range("E19").Value = 0.99999770911318

# Row 20 (C20 flips true -> false)
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = [double]"4.771315079648401E-10"
$ws.Range("E20").Value = 0.9999999995228684

# Row 21 (C21 flips true -> false)
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = 0.09296938557109211
$ws.Range("E21").Value = 0.9070306144289079
$ws.Range("F21").Value = 11.0827465057373
$ws.Range("G21").Value = 0.2
